$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 23.34546866666667
$ws.Cells.Item(2, 8).Value = 70.036406
$ws.Cells.Item(2, 9).Value = 0.4715073400272545
$ws.Cells.Item(2, 10).Value = 0.4715073400272545
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 97.15988866666665
$ws.Cells.Item(2, 14).Value = 291.479666
$ws.Cells.Item(2, 15).Value = 0.974735784617843
$ws.Cells.Item(2, 16).Value = 0.9747357846178432
$ws.Cells.Item(2, 17).Value = 2268.243136524488
$ws.Cells.Item(2, 18).Value = 20414.18822872039
$ws.Cells.Item(2, 19).Value = 0.459595077034538
$ws.Cells.Item(2, 20).Value = 0.4595950770345381
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 23.34546866666667
$ws.Cells.Item(3, 8).Value = 70.036406
$ws.Cells.Item(3, 9).Value = 0.4715073400272545
$ws.Cells.Item(3, 10).Value = 0.4715073400272545
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.781603
$ws.Cells.Item(3, 14).Value = 2.344809
$ws.Cells.Item(3, 15).Value = 0.007841264784466923
$ws.Cells.Item(3, 16).Value = 0.007841264784466923
$ws.Cells.Item(3, 17).Value = 18.24688834627267
$ws.Cells.Item(3, 18).Value = 164.221995116454
$ws.Cells.Item(3, 19).Value = 0.003697213900973381
$ws.Cells.Item(3, 20).Value = 0.003697213900973381
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 23.34546866666667
$ws.Cells.Item(4, 8).Value = 70.036406
$ws.Cells.Item(4, 9).Value = 0.4715073400272545
$ws.Cells.Item(4, 10).Value = 0.4715073400272545
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.736688
$ws.Cells.Item(4, 14).Value = 5.210064
$ws.Cells.Item(4, 15).Value = 0.01742295059768999
$ws.Cells.Item(4, 16).Value = 0.01742295059768999
$ws.Cells.Item(4, 17).Value = 40.543795287776
$ws.Cells.Item(4, 18).Value = 364.894157589984
$ws.Cells.Item(4, 19).Value = 0.008215049091743072
$ws.Cells.Item(4, 20).Value = 0.008215049091743072
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.57337833333333
$ws.Cells.Item(5, 8).Value = 37.720135
$ws.Cells.Item(5, 9).Value = 0.2539439348061199
$ws.Cells.Item(5, 10).Value = 0.2539439348061199
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 97.15988866666665
$ws.Cells.Item(5, 14).Value = 291.479666
$ws.Cells.Item(5, 15).Value = 0.974735784617843
$ws.Cells.Item(5, 16).Value = 0.9747357846178432
$ws.Cells.Item(5, 17).Value = 1221.628039030545
$ws.Cells.Item(5, 18).Value = 10994.65235127491
$ws.Cells.Item(5, 19).Value = 0.2475282405421856
$ws.Cells.Item(5, 20).Value = 0.2475282405421857
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.57337833333333
$ws.Cells.Item(6, 8).Value = 37.720135
$ws.Cells.Item(6, 9).Value = 0.2539439348061199
$ws.Cells.Item(6, 10).Value = 0.2539439348061199
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.781603
$ws.Cells.Item(6, 14).Value = 2.344809
$ws.Cells.Item(6, 15).Value = 0.007841264784466923
$ws.Cells.Item(6, 16).Value = 0.007841264784466923
$ws.Cells.Item(6, 17).Value = 9.827390225468333
$ws.Cells.Item(6, 18).Value = 88.446512029215
$ws.Cells.Item(6, 19).Value = 0.001991241633224192
$ws.Cells.Item(6, 20).Value = 0.001991241633224192
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.57337833333333
$ws.Cells.Item(7, 8).Value = 37.720135
$ws.Cells.Item(7, 9).Value = 0.2539439348061199
$ws.Cells.Item(7, 10).Value = 0.2539439348061199
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.736688
$ws.Cells.Item(7, 14).Value = 5.210064
$ws.Cells.Item(7, 15).Value = 0.01742295059768999
$ws.Cells.Item(7, 16).Value = 0.01742295059768999
$ws.Cells.Item(7, 17).Value = 21.83603527096
$ws.Cells.Item(7, 18).Value = 196.52431743864
$ws.Cells.Item(7, 19).Value = 0.004424452630710034
$ws.Cells.Item(7, 20).Value = 0.004424452630710034
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 13.59357133333333
$ws.Cells.Item(8, 8).Value = 40.780714
$ws.Cells.Item(8, 9).Value = 0.2745487251666257
$ws.Cells.Item(8, 10).Value = 0.2745487251666257
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 97.15988866666665
$ws.Cells.Item(8, 14).Value = 291.479666
$ws.Cells.Item(8, 15).Value = 0.974735784617843
$ws.Cells.Item(8, 16).Value = 0.9747357846178432
$ws.Cells.Item(8, 17).Value = 1320.749877329058
$ws.Cells.Item(8, 18).Value = 11886.74889596152
$ws.Cells.Item(8, 19).Value = 0.2676124670411195
$ws.Cells.Item(8, 20).Value = 0.2676124670411195
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 13.59357133333333
$ws.Cells.Item(9, 8).Value = 40.780714
$ws.Cells.Item(9, 9).Value = 0.2745487251666257
$ws.Cells.Item(9, 10).Value = 0.2745487251666257
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.781603
$ws.Cells.Item(9, 14).Value = 2.344809
$ws.Cells.Item(9, 15).Value = 0.007841264784466923
$ws.Cells.Item(9, 16).Value = 0.007841264784466923
$ws.Cells.Item(9, 17).Value = 10.62477613484733
$ws.Cells.Item(9, 18).Value = 95.622985213626
$ws.Cells.Item(9, 19).Value = 0.00215280925026935
$ws.Cells.Item(9, 20).Value = 0.00215280925026935
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 13.59357133333333
$ws.Cells.Item(10, 8).Value = 40.780714
$ws.Cells.Item(10, 9).Value = 0.2745487251666257
$ws.Cells.Item(10, 10).Value = 0.2745487251666257
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 1.736688
$ws.Cells.Item(10, 14).Value = 5.210064
$ws.Cells.Item(10, 15).Value = 0.01742295059768999
$ws.Cells.Item(10, 16).Value = 0.01742295059768999
$ws.Cells.Item(10, 17).Value = 23.607792211744
$ws.Cells.Item(10, 18).Value = 212.470129905696
$ws.Cells.Item(10, 19).Value = 0.004783448875236886
$ws.Cells.Item(10, 20).Value = 0.004783448875236886
